$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds numeric-looking text (e.g. "508.99", "1.00") that
# must stay stored as text, not auto-converted to a number. Force the cell
# to Text format before writing, then restore the default style so no
# stray formatting is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.359.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.324.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "508.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0995"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.69%  "
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.336"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.740.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.320.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.330.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "321.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +10.13%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.162"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.86%  "
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0714"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.877"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "149.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.36%  "
$ws.Range("E41").Value = "  -1.79%  "
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "276.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0923"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0493"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "17.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.376"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0212"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.31%  "
